$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.881.40"
$ws.Range("E2").Value = "  +1.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.966.61"
$ws.Range("E3").Value = "  +4.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9912"
$ws.Range("E4").Value = "  -0.86%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8334"
$ws.Range("E5").Value = "  +76.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "252.42"
$ws.Range("E6").Value = "  +3.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9866"
$ws.Range("E7").Value = "  -1.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3460"
$ws.Range("E8").Value = "  +19.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.55"
$ws.Range("E9").Value = "  +14.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06864"
$ws.Range("E10").Value = "  +5.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8597"
$ws.Range("E11").Value = "  +18.26%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08068"
$ws.Range("E12").Value = "  +3.98%  "

$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "101.38"
$ws.Range("E13").Value = "  +5.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.961.50"
$ws.Range("E14").Value = "  +3.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.441"
$ws.Range("E15").Value = "  +4.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.60"
$ws.Range("E16").Value = "  -1.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.815.89"
$ws.Range("E17").Value = "  +1.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.82"
$ws.Range("E18").Value = "  +5.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007814"
$ws.Range("E19").Value = "  +4.46%  "

$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.225.72"
$ws.Range("E20").Value = "  +4.17%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.645"
$ws.Range("E21").Value = "  +6.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9838"
$ws.Range("E22").Value = "  -1.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9780"
$ws.Range("E23").Value = "  -2.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.732"
$ws.Range("E24").Value = "  +6.51%  "

$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1698"
$ws.Range("E25").Value = "  +74.98%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.578"
$ws.Range("E26").Value = "  +5.21%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.97"
$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.45"
$ws.Range("E28").Value = "  +2.95%  "

$ws.Range("E29").Value = "  +15.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.555"
$ws.Range("E30").Value = "  +5.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.345"
$ws.Range("E31").Value = "  +0.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.504"
$ws.Range("E32").Value = "  +5.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.292"
$ws.Range("E33").Value = "  +3.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05082"
$ws.Range("E34").Value = "  +4.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.215"
$ws.Range("E35").Value = "  +7.68%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7376"
$ws.Range("E36").Value = "  +6.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.799"
$ws.Range("E37").Value = "  +3.07%  "

$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9802"
$ws.Range("E38").Value = "  -1.91%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01978"
$ws.Range("E39").Value = "  +4.63%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.870"
$ws.Range("E40").Value = "  +1.94%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.549"
$ws.Range("E41").Value = "  +5.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "77.67"
$ws.Range("E42").Value = "  +3.21%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4633"
$ws.Range("E43").Value = "  +8.30%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.065"
$ws.Range("E44").Value = "  +4.01%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8479"
$ws.Range("E45").Value = "  +2.17%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.71"
$ws.Range("E46").Value = "  +2.28%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9836"
$ws.Range("E47").Value = "  -1.66%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.925"
$ws.Range("E48").Value = "  +2.89%  "

$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.446"
$ws.Range("E49").Value = "  +6.85%  "

$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4262"
$ws.Range("E50").Value = "  +7.74%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.18"
$ws.Range("E51").Value = "  +2.71%  "
